$d = $word.ActiveDocument

# 1) Replace "line-height: 3;" with fullwidth-3 version
$d.Content.Find.Execute("line-height: 3;", $true, $false, $false, $false, $false,
                         $true, 1, $false, "line-height: ３;", 2)

# 2) Replace "font-size: 1.5em;" with fullwidth version
$d.Content.Find.Execute("font-size: 1.5em;", $true, $false, $false, $false, $false,
                         $true, 1, $false, "font-size: １.５em;", 2)

# 3) Find the bullet-list paragraph ("• " followed by "1") and restructure it
foreach ($p in $d.Paragraphs) {
    $ptext = $p.Range.Text
    if ($ptext -like "*$([char]0x2022)*") {
        # This is the bullet paragraph.
        $p.Range.Text = [string][char]0xFF11
        $p.Style = "ListBullet"
        $p.Range.Font.NameAscii = "台灣明體"
        $p.Range.Font.NameFarEast = "台灣明體"
        $p.Range.Font.Name = "台灣明體"
        $p.Range.Font.Size = 13
        break
    }
}
